# Auto-sync batch FINAL (AUTO-TIMEOUT)
# Appends 6 new reconciliation rows (89-94) to Sheet1, extending the
# used range from A1:I88 to A1:I94.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ A = 237674890585; B = "LA NEGRESSE LTDLA CBOX R0 NGO MBOCK epse MBAYAN MARIE CLAIRE"; C = 0;        D = "Pk8";                        E = 10000;    F = 67456;  G = 57456;   H = 6.7456;               I = "Ndogbong" },
    @{ A = 237653294562; B = "NANHOU KEMAYOU AVIGAEL ETS MOBILE FINANCIAL SERVICES MFS";       C = "Rte_5"; D = "Socaver Ndongbong";          E = 162167.4; F = 506519; G = 344351.6; H = 3.123432946449163;    I = "Ndogbong" },
    @{ A = 237681659043; B = "SYLVIE DJIDJOU TEGUIA EPSE TOUKOU";                               C = 0;        D = "Makepe Conquete";            E = 10000;    F = 18480;  G = 8480;    H = 1.848;                I = "Ndogbong" },
    @{ A = 237679086144; B = "ALAIN CHETEU KAMDEM";                                             C = 0;        D = "Cite Sic Stade Marion Ocm";  E = 10000;    F = 2615;   G = -7385;   H = 0.2615;               I = "Cite Sic" },
    @{ A = 237683998069; B = "MEDJEU FEUZEU epse FEGHEM WAHOUE JOSIANE ETS MOBILE FINANCIAL SERVICES MFS"; C = 0; D = "Makepe Conquete";      E = 20800;    F = 17949;  G = -2851;   H = 0.8629326923076923;   I = "Ndogbong" },
    @{ A = 237679252522; B = "WARAMMA NICOLAS";                                                 C = 0;        D = "Cite Sic Stade Marion Ocm";  E = 10000;    F = 23;     G = -9977;   H = 0.0023;               I = "Cite Sic" }
)

$startRow = 89
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value2 = $row.A
    $ws.Cells.Item($r, 2).Value2 = $row.B
    $ws.Cells.Item($r, 3).Value2 = $row.C
    $ws.Cells.Item($r, 4).Value2 = $row.D
    $ws.Cells.Item($r, 5).Value2 = $row.E
    $ws.Cells.Item($r, 6).Value2 = $row.F
    $ws.Cells.Item($r, 7).Value2 = $row.G
    $ws.Cells.Item($r, 8).Value2 = $row.H
    $ws.Cells.Item($r, 9).Value2 = $row.I
}
